$d = $word.ActiveDocument
$count = $d.Paragraphs.Count
Write-Host "Paragraph count:" $count
for ($i = 1; $i -le $count; $i++) {
  $p = $d.Paragraphs.Item($i)
  $rng = $p.Range
  $pxml = $rng.WordOpenXML
  if ($pxml -match '<w:contextualSpacing[^/]*/>') {
    $newxml = $pxml -replace '<w:contextualSpacing[^/]*/>', ''
    $rng.InsertXML($newxml)
  }
}
Write-Host "Done all paragraphs"
